# Update cryptocurrency price/volume figures per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.908.73"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "1.639.25"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'213.69"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'23.67"
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").Value = "'0.0875"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "1.640.04"
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").Value = "'0.574"
$ws.Range("E15").Value = "  +4.30%  "
$ws.Range("D16").Value = "'66.19"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").Value = "27.897.99"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").Value = "'232.06"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "'7.60"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'10.94"
$ws.Range("E22").Value = "  +4.61%  "
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").Value = "'2.07"
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("D25").Value = "'151.75"
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("D26").Value = "'6.92"
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").Value = "'15.71"
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("D33").Value = "'3.12"
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("D34").Value = "1.418.44"
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("E37").Value = "  +2.24%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'0.556"
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.918"
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("E41").Value = "  +1.72%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "'67.09"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("E45").Value = "  +3.67%  "
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").Value = "1.780.63"
$ws.Range("D49").Value = "'88.53"
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("E51").Value = "  +0.65%  "
